$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-31 Tuesday" "2025-01-01 Wednesday"

Replace-Text "991×5=" "921×7="
Replace-Text "518×8=" "199×5="
Replace-Text "923×7=" "119×3="
Replace-Text "395×8=" "780×3="
Replace-Text "384×8=" "830×2="

Replace-Text "355×9=" "451×8="
Replace-Text "288×9=" "363×4="
Replace-Text "750×2=" "186×9="
Replace-Text "637×7=" "402×2="
Replace-Text "320×9=" "322×6="

Replace-Text "523×4=" "814×4="
Replace-Text "336×8=" "742×2="
Replace-Text "914×2=" "785×5="
Replace-Text "520×6=" "103×4="
Replace-Text "862×4=" "461×3="

Replace-Text "821×2=" "764×2="
Replace-Text "291×3=" "970×2="
Replace-Text "804×9=" "367×8="
Replace-Text "351×2=" "489×4="
Replace-Text "392×8=" "435×4="

Replace-Text "425×7=" "673×6="
Replace-Text "835×2=" "687×8="
Replace-Text "585×9=" "788×4="
Replace-Text "776×5=" "363×8="
Replace-Text "476×6=" "690×7="

"Done"
